# Add a new blank slide at the end of the presentation (slide 5),
# mirroring "Random Button Added, AutoPlay Button Added" groundwork:
# a fresh blank slide using the "Leer" (Blank) layout.

$p = $ppt.ActivePresentation

# Append a new slide after the last existing slide, using the blank layout
# (ppLayoutBlank = 12), matching PowerPoint's "New Slide" behaviour.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)
